# Auto-generated Excel COM-interop script
# Updates column F ('想去人数' / interested-attendee counts) values
# on worksheets '展览' (sheet 1) and '全部类型' (sheet 4)
# per the commit 'Update gh-pages to output generated at 456a3b4'.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # 展览
$ws1.Range("F3").Value = 1180
$ws1.Range("F4").Value = 1611
$ws1.Range("F5").Value = 186
$ws1.Range("F6").Value = 186
$ws1.Range("F8").Value = 1568
$ws1.Range("F9").Value = 3177
$ws1.Range("F10").Value = 722
$ws1.Range("F11").Value = 1904
$ws1.Range("F12").Value = 1856
$ws1.Range("F13").Value = 930
$ws1.Range("F14").Value = 319
$ws1.Range("F16").Value = 1542
$ws1.Range("F17").Value = 313
$ws1.Range("F19").Value = 46
$ws1.Range("F20").Value = 1337
$ws1.Range("F21").Value = 449
$ws1.Range("F22").Value = 549
$ws1.Range("F23").Value = 238
$ws1.Range("F24").Value = 7941
$ws1.Range("F25").Value = 9286
$ws1.Range("F26").Value = 799
$ws1.Range("F27").Value = 611
$ws1.Range("F28").Value = 1760
$ws1.Range("F29").Value = 115
$ws1.Range("F30").Value = 301

$ws4 = $wb.Worksheets.Item(4)   # 全部类型
$ws4.Range("F4").Value = 1180
$ws4.Range("F5").Value = 1611
$ws4.Range("F6").Value = 186
$ws4.Range("F7").Value = 186
$ws4.Range("F10").Value = 1568
$ws4.Range("F11").Value = 3177
$ws4.Range("F12").Value = 722
$ws4.Range("F13").Value = 1904
$ws4.Range("F14").Value = 1856
$ws4.Range("F15").Value = 930
$ws4.Range("F16").Value = 319
$ws4.Range("F18").Value = 1542
$ws4.Range("F19").Value = 313
$ws4.Range("F22").Value = 46
$ws4.Range("F24").Value = 1337
$ws4.Range("F25").Value = 449
$ws4.Range("F26").Value = 549
$ws4.Range("F27").Value = 238
$ws4.Range("F28").Value = 7941
$ws4.Range("F29").Value = 9286
$ws4.Range("F30").Value = 799
$ws4.Range("F31").Value = 611
$ws4.Range("F32").Value = 1760
$ws4.Range("F35").Value = 115
$ws4.Range("F36").Value = 301
